$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2025/12/03 02:00"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "-"
$ws.Range("G11").Value = "-"
